# Updated via Streamlit Approval System
# Applies the pending-approval sheet edits described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextZero([string]$addr) {
    # Force the literal text "0" (not the number 0) into the cell, then
    # drop back to the default style so no stray NumberFormat lingers.
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = "0"
    $c.Style = "Normal"
}

# --- Row 12 ---
$ws.Range("V12").Value = 7000

$ws.Range("AI12").Value = ""
$ws.Range("AJ12").Value = ""

Set-TextZero "AK12"
Set-TextZero "AL12"
Set-TextZero "AM12"
Set-TextZero "AN12"
Set-TextZero "AO12"

# --- Row 13 ---
Set-TextZero "AK13"
Set-TextZero "AL13"
Set-TextZero "AM13"
Set-TextZero "AN13"
Set-TextZero "AO13"

# --- Row 18 ---
Set-TextZero "AN18"

# --- Row 19 ---
Set-TextZero "AN19"

# --- Row 20 ---
Set-TextZero "AN20"

# --- Row 21 ---
Set-TextZero "AN21"

# --- Row 22 ---
Set-TextZero "AN22"

# --- Row 23 ---
Set-TextZero "AN23"

# --- Row 24 ---
Set-TextZero "AN24"

# --- Row 25 ---
Set-TextZero "AK25"

# --- Row 26 ---
Set-TextZero "AK26"

# --- Row 27 ---
Set-TextZero "AK27"
Set-TextZero "AL27"
Set-TextZero "AM27"
Set-TextZero "AN27"
Set-TextZero "AO27"
